# Apply commit:
#  "update strings to include the download link on the start screen
#   move icon to its own folder
#   fade out settings button to make it less imposing"
#
# The part of this change that touches the workbook is the new
# START_SCREEN_INSPECT_TEXT / download-link string that was added as a new
# row in the "BasicText" localization sheet (row 49), pushing the existing
# tutorial-text rows down by one. Two sheet selections also moved as a
# side-effect of editing in Excel.

$wb = $excel.ActiveWorkbook

# --- StringLocalizations_BasicText: insert the new localization row ---
$wsBasic = $wb.Worksheets.Item("StringLocalizations_BasicText")

# Shift existing rows 49:56 down to 50:57 and create a fresh row 49.
$wsBasic.Rows(49).Insert()

$wsBasic.Range("A49").Value = "START_SCREEN_INSPECT_TEXT"
$wsBasic.Range("B49").Value = "Help your local police force to solve real crimes with the Inspec2t app"
$wsBasic.Range("C49").Value = "XXXX"
$wsBasic.Range("D49").Value = "XXXX"
$wsBasic.Range("E49").Value = "XXXX"

# Reflect where the editor ended up after making the change.
$wsBasic.Activate()
$wsBasic.Range("A50").Select()

# --- StringLocalizations_Valencia: selection moved while reviewing ---
$wsValencia = $wb.Worksheets.Item("StringLocalizations_Valencia")
$wsValencia.Activate()
$wsValencia.Range("E5").Select()

# Leave the BasicText sheet as the active/selected tab, matching the
# workbook's original active tab.
$wsBasic.Activate()
